# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 644
$wsExhibit.Range("F4").Value = 209
$wsExhibit.Range("F6").Value = 9764
$wsExhibit.Range("F10").Value = 3347
$wsExhibit.Range("F15").Value = 282
$wsExhibit.Range("F16").Value = 527
$wsExhibit.Range("F17").Value = 106
$wsExhibit.Range("F19").Value = 1422

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 644
$wsAll.Range("F5").Value = 209
$wsAll.Range("F7").Value = 9764
$wsAll.Range("F11").Value = 3347
$wsAll.Range("F16").Value = 282
$wsAll.Range("F17").Value = 527
$wsAll.Range("F18").Value = 106
$wsAll.Range("F20").Value = 1422
